# Auto-generated script to update worksheet cell values per scraped-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 730.2069  # H17: was 810.2963
$ws.Cells.Item(17, 10).Value = 769.4815  # J17: was 859.12
$ws.Cells.Item(17, 12).Value = 2308.4445  # L17: was 2577.36
$ws.Cells.Item(17, 14).Value = -2644.4445  # N17: was -2913.36
$ws.Cells.Item(98, 8).Value = 1350.1875  # H98: was 1350.8125
$ws.Cells.Item(98, 9).Value = 1050.3572  # I98: was 1051.0714
$ws.Cells.Item(98, 11).Value = 1050.3572  # K98: was 1051.0714
$ws.Cells.Item(98, 13).Value = 447.6428000000001  # M98: was 446.9286
$ws.Cells.Item(111, 8).Value = 1287.1111  # H111: was 1361
$ws.Cells.Item(111, 10).Value = 697  # J111: was 697.5
$ws.Cells.Item(111, 12).Value = 2091  # L111: was 2092.5
$ws.Cells.Item(111, 14).Value = -8225  # N111: was -8226.5
$ws.Cells.Item(113, 8).Value = 7299.5  # H113: was 7232.6665
$ws.Cells.Item(113, 9).Value = 6000  # I113: was 6549.5
$ws.Cells.Item(113, 11).Value = 6000  # K113: was 6549.5
$ws.Cells.Item(113, 13).Value = -2746  # M113: was -3295.5
$ws.Cells.Item(122, 8).Value = 1350.1875  # H122: was 1350.8125
$ws.Cells.Item(122, 9).Value = 1050.3572  # I122: was 1051.0714
$ws.Cells.Item(122, 11).Value = 3151.0716  # K122: was 3153.2142
$ws.Cells.Item(122, 13).Value = -701.0715999999998  # M122: was -703.2142000000003
$ws.Cells.Item(137, 8).Value = 2328.077  # H137: was 2164.7334
$ws.Cells.Item(137, 9).Value = 1900.125  # I137: was 1785.5555
$ws.Cells.Item(137, 10).Value = 3012.8  # J137: was 2733.5
$ws.Cells.Item(137, 11).Value = 5700.375  # K137: was 5356.666499999999
$ws.Cells.Item(137, 12).Value = 9038.400000000001  # L137: was 8200.5
$ws.Cells.Item(137, 13).Value = -3150.375  # M137: was -2806.666499999999
$ws.Cells.Item(137, 14).Value = -14138.4  # N137: was -13300.5
$ws.Cells.Item(141, 8).Value = 2874.697  # H141: was 3137.2058
$ws.Cells.Item(141, 9).Value = 2683.2812  # I141: was 2770.8125
$ws.Cells.Item(141, 10).Value = 9000  # J141: was 8999.5
$ws.Cells.Item(141, 11).Value = 8049.8436  # K141: was 8312.4375
$ws.Cells.Item(141, 12).Value = 27000  # L141: was 26998.5
$ws.Cells.Item(141, 13).Value = -2869.8436  # M141: was -3132.4375
$ws.Cells.Item(141, 14).Value = -37360  # N141: was -37358.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(12, 8).Value = 650  # H12: was 700
$ws.Cells.Item(12, 9).Value = 650  # I12: was 0
$ws.Cells.Item(12, 10).Value = 0  # J12: was 700
$ws.Cells.Item(12, 11).Value = 650  # K12: was 0
$ws.Cells.Item(12, 12).Value = 0  # L12: was 700
$ws.Cells.Item(12, 13).Value = -477  # M12: was None
$ws.Cells.Item(12, 14).ClearContents()  # N12: was -1046
$ws.Cells.Item(122, 8).Value = 1566.5294  # H122: was 1684.3889
$ws.Cells.Item(122, 9).Value = 1545.6875  # I122: was 1655.9333
$ws.Cells.Item(122, 10).Value = 1900  # J122: was 1826.6666
$ws.Cells.Item(122, 11).Value = 4637.0625  # K122: was 4967.7999
$ws.Cells.Item(122, 12).Value = 5700  # L122: was 5479.9998
$ws.Cells.Item(122, 13).Value = -2187.0625  # M122: was -2517.7999
$ws.Cells.Item(122, 14).Value = -10600  # N122: was -10379.9998
$ws.Cells.Item(135, 8).Value = 44999.5  # H135: was 45000
$ws.Cells.Item(135, 10).Value = 44999.5  # J135: was 45000
$ws.Cells.Item(135, 12).Value = 44999.5  # L135: was 45000
$ws.Cells.Item(135, 14).Value = -55139.5  # N135: was -55140

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2153.375  # H86: was 2032.4286
$ws.Cells.Item(86, 10).Value = 3000  # J86: was 0
$ws.Cells.Item(86, 12).Value = 3000  # L86: was 0
$ws.Cells.Item(86, 14).Value = -5246  # N86: was None
$ws.Cells.Item(89, 8).Value = 2153.375  # H89: was 2032.4286
$ws.Cells.Item(89, 10).Value = 3000  # J89: was 0
$ws.Cells.Item(89, 12).Value = 15000  # L89: was 0
$ws.Cells.Item(89, 14).Value = -26232  # N89: was None
$ws.Cells.Item(135, 8).Value = 59997  # H135: was 59999
$ws.Cells.Item(135, 10).Value = 59997  # J135: was 59999
$ws.Cells.Item(135, 12).Value = 59997  # L135: was 59999
$ws.Cells.Item(135, 14).Value = -70137  # N135: was -70139

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(69, 8).Value = 40199.332  # H69: was 40200
$ws.Cells.Item(69, 10).Value = 40199.332  # J69: was 40200
$ws.Cells.Item(69, 12).Value = 40199.332  # L69: was 40200
$ws.Cells.Item(69, 14).Value = -41697.332  # N69: was -41698
$ws.Cells.Item(72, 8).Value = 40199.332  # H72: was 40200
$ws.Cells.Item(72, 10).Value = 40199.332  # J72: was 40200
$ws.Cells.Item(72, 12).Value = 120597.996  # L72: was 120600
$ws.Cells.Item(72, 14).Value = -128085.996  # N72: was -128088
$ws.Cells.Item(93, 8).Value = 26125  # H93: was 23788.2
$ws.Cells.Item(93, 9).Value = 12000  # I93: was 10000
$ws.Cells.Item(93, 10).Value = 28950  # J93: was 27235.25
$ws.Cells.Item(93, 11).Value = 12000  # K93: was 10000
$ws.Cells.Item(93, 12).Value = 28950  # L93: was 27235.25
$ws.Cells.Item(93, 13).Value = -10128  # M93: was -8128
$ws.Cells.Item(93, 14).Value = -32694  # N93: was -30979.25
$ws.Cells.Item(107, 8).Value = 3902.6667  # H107: was 4207.091
$ws.Cells.Item(107, 9).Value = 847.75  # I107: was 889.7143
$ws.Cells.Item(107, 11).Value = 847.75  # K107: was 889.7143
$ws.Cells.Item(107, 13).Value = 1072.25  # M107: was 1030.2857
$ws.Cells.Item(122, 8).Value = 51940  # H122: was 24344.455
$ws.Cells.Item(122, 9).Value = 0  # I122: was 1137.8
$ws.Cells.Item(122, 10).Value = 51940  # J122: was 43683.332
$ws.Cells.Item(122, 11).Value = 0  # K122: was 3413.4
$ws.Cells.Item(122, 12).Value = 155820  # L122: was 131049.996
$ws.Cells.Item(122, 13).ClearContents()  # M122: was -963.3999999999996
$ws.Cells.Item(122, 14).Value = -160720  # N122: was -135949.996

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(63, 8).Value = 500  # H63: was 0
$ws.Cells.Item(63, 9).Value = 500  # I63: was 0
$ws.Cells.Item(63, 11).Value = 1500  # K63: was 0
$ws.Cells.Item(63, 13).Value = -751  # M63: was None
$ws.Cells.Item(66, 8).Value = 500  # H66: was 0
$ws.Cells.Item(66, 9).Value = 500  # I66: was 0
$ws.Cells.Item(66, 11).Value = 4500  # K66: was 0
$ws.Cells.Item(66, 13).Value = -756  # M66: was None
$ws.Cells.Item(113, 8).Value = 453.07144  # H113: was 419.23077
$ws.Cells.Item(113, 10).Value = 897.2  # J113: was 898.25
$ws.Cells.Item(113, 12).Value = 2691.6  # L113: was 2694.75
$ws.Cells.Item(113, 14).Value = -7031.6  # N113: was -7034.75
$ws.Cells.Item(134, 8).Value = 14762.174  # H134: was 14224.167
$ws.Cells.Item(134, 9).Value = 13176.667  # I134: was 10345
$ws.Cells.Item(134, 11).Value = 39530.001  # K134: was 31035
$ws.Cells.Item(134, 13).Value = -34460.001  # M134: was -25965

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3174  # H102: was 2569.25
$ws.Cells.Item(102, 9).Value = 2484.2  # I102: was 1649.6666
$ws.Cells.Item(102, 10).Value = 3439.3076  # J102: was 3121
$ws.Cells.Item(102, 11).Value = 2484.2  # K102: was 1649.6666
$ws.Cells.Item(102, 12).Value = 3439.3076  # L102: was 3121
$ws.Cells.Item(102, 13).Value = -862.1999999999998  # M102: was -27.66660000000002
$ws.Cells.Item(102, 14).Value = -6683.3076  # N102: was -6365
$ws.Cells.Item(113, 8).Value = 2999  # H113: was 2095
$ws.Cells.Item(113, 9).Value = 2999  # I113: was 2190
$ws.Cells.Item(113, 10).Value = 0  # J113: was 2000
$ws.Cells.Item(113, 11).Value = 2999  # K113: was 2190
$ws.Cells.Item(113, 12).Value = 0  # L113: was 2000
$ws.Cells.Item(113, 13).Value = -829  # M113: was -20
$ws.Cells.Item(113, 14).ClearContents()  # N113: was -6340
$ws.Cells.Item(122, 8).Value = 4022.8572  # H122: was 4326.75
$ws.Cells.Item(122, 9).Value = 4630.5  # I122: was 4900.6665
$ws.Cells.Item(122, 10).Value = 2503.75  # J122: was 2605
$ws.Cells.Item(122, 11).Value = 13891.5  # K122: was 14701.9995
$ws.Cells.Item(122, 12).Value = 7511.25  # L122: was 7815
$ws.Cells.Item(122, 13).Value = -11441.5  # M122: was -12251.9995
$ws.Cells.Item(122, 14).Value = -12411.25  # N122: was -12715
$ws.Cells.Item(126, 8).Value = 2785  # H126: was 2900
$ws.Cells.Item(126, 9).Value = 2440  # I126: was 0
$ws.Cells.Item(126, 11).Value = 7320  # K126: was 0
$ws.Cells.Item(126, 13).Value = -4850  # M126: was None

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3885.2222  # H16: was 2120.875
$ws.Cells.Item(16, 10).Value = 18000  # J16: was 0
$ws.Cells.Item(16, 12).Value = 18000  # L16: was 0
$ws.Cells.Item(16, 14).Value = -18340  # N16: was None
$ws.Cells.Item(62, 8).Value = 49374.5  # H62: was 60248.5
$ws.Cells.Item(62, 10).Value = 49374.5  # J62: was 60248.5
$ws.Cells.Item(62, 12).Value = 49374.5  # L62: was 60248.5
$ws.Cells.Item(62, 14).Value = -50622.5  # N62: was -61496.5
$ws.Cells.Item(65, 8).Value = 49374.5  # H65: was 60248.5
$ws.Cells.Item(65, 10).Value = 49374.5  # J65: was 60248.5
$ws.Cells.Item(65, 12).Value = 148123.5  # L65: was 180745.5
$ws.Cells.Item(65, 14).Value = -154363.5  # N65: was -186985.5
$ws.Cells.Item(68, 8).Value = 5694.65  # H68: was 5841.737
$ws.Cells.Item(68, 9).Value = 4192.533  # I68: was 4199.2
$ws.Cells.Item(68, 10).Value = 10201  # J68: was 12001.25
$ws.Cells.Item(68, 11).Value = 4192.533  # K68: was 4199.2
$ws.Cells.Item(68, 12).Value = 10201  # L68: was 12001.25
$ws.Cells.Item(68, 13).Value = -3443.533  # M68: was -3450.2
$ws.Cells.Item(68, 14).Value = -11699  # N68: was -13499.25
$ws.Cells.Item(71, 8).Value = 5694.65  # H71: was 5841.737
$ws.Cells.Item(71, 9).Value = 4192.533  # I71: was 4199.2
$ws.Cells.Item(71, 10).Value = 10201  # J71: was 12001.25
$ws.Cells.Item(71, 11).Value = 20962.665  # K71: was 20996
$ws.Cells.Item(71, 12).Value = 51005  # L71: was 60006.25
$ws.Cells.Item(71, 13).Value = -17218.665  # M71: was -17252
$ws.Cells.Item(71, 14).Value = -58493  # N71: was -67494.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1683.6  # H100: was 1804
$ws.Cells.Item(100, 9).Value = 978.2857  # I100: was 1041.3334
$ws.Cells.Item(100, 11).Value = 1956.5714  # K100: was 2082.6668
$ws.Cells.Item(100, 13).Value = -1415.5714  # M100: was -1541.6668
$ws.Cells.Item(122, 8).Value = 3043.7334  # H122: was 3046.4
$ws.Cells.Item(122, 9).Value = 3400.6  # I122: was 3404.6
$ws.Cells.Item(122, 11).Value = 10201.8  # K122: was 10213.8
$ws.Cells.Item(122, 13).Value = -7751.799999999999  # M122: was -7763.799999999999
